# Edit LOQ4217.xlsx worksheet: insert a new row (Docentes responsaveis data row)
# and update several rows' B/C content (Objetivos, Programa resumido, Programa,
# Metodo/Criterio/Norma de recuperacao/Bibliografia) to reflect the new long-form
# text content, shifting everything below row 12 down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row before row 13. This shifts old rows 13..23 down to 14..24
#    and conveniently carries the row-height formatting down with them, which
#    already matches the target layout for every row except the (new) row 15.
$ws.Rows.Item(13).Insert()

# The insert leaves a stray formatted-but-empty A13 cell behind (inherited
# style from the row above) - remove it completely so the row has no cells,
# matching the target (row 13 only has B13/C13).
$ws.Range("A13").Clear()

# 2. Fix up the one row whose custom height should NOT carry over: row 15
#    ("Short syllabus:") - after the shift it kept the 60pt custom height that
#    belonged to the old "Objetivos:" row, but the target has no custom height
#    for this row, so auto-fit it back to the default.
$ws.Rows.Item(15).AutoFit()

# 3. Populate the new row 13 with the "Docentes responsaveis" value that used
#    to sit in old row 13 (B/C), now moved here. The inserted row's B/C cells
#    do not exist yet, so first clone the B/C formatting (wrap text, vertical
#    top alignment, black/red font) from an existing untouched data row (row
#    3) before writing the values, so the same style indices (s="2"/s="3")
#    are reused instead of new style entries being created.
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

# 4. Row 10 ("Objetivos:") B/C - new long description replacing the old
#    "5840560 - Marco Antonio Carvalho Pereira" placeholder text.
$ws.Range("B10").Value = "Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos."
$ws.Range("C10").Value = "Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos."

# 5. Row 14 ("Programa resumido:") B/C - replaces "Semestral" with the short
#    syllabus outline text.
$ws.Range("B14").Value = "1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa"
$ws.Range("C14").Value = "1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa"

# 6. Row 16 ("Programa:") B/C - replaces "01/01/2021" with the full syllabus
#    text.
$ws.Range("B16").Value = "1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos."
$ws.Range("C16").Value = "1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos."

# 7. Row 19 ("Método:") B/C - now holds the "Aulas expositivas..." text
#    (previously on what is now row 20).
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas de exercícios."

# 8. Row 20 ("Critério:") B/C - now holds the "Média de Provas..." text
#    (previously on what is now row 21).
$ws.Range("B20").Value = "Média de Provas e trabalhos (MF)."
$ws.Range("C20").Value = "Média de Provas e trabalhos (MF)."

# 9. Row 21 ("Norma de recuperação:") B/C - now holds the "Prova de
#    Recuperação..." text (previously on what is now row 22).
$ws.Range("B21").Value = "Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR"
$ws.Range("C21").Value = "Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR"

# 10. Row 22 ("Bibliografia:") B/C - new bibliography text.
$ws.Range("B22").Value = "CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012"
$ws.Range("C22").Value = "CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012"
